$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.02514266666667
$ws.Range("H2").Value = 45.075428
$ws.Range("I2").Value = 0.1401726531301337
$ws.Range("J2").Value = 0.1401726531301337
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 68.00339
$ws.Range("N2").Value = 204.01017
$ws.Range("O2").Value = 0.6265962299909886
$ws.Range("P2").Value = 0.6265962299909885
$ws.Range("Q2").Value = 1021.760636566973
$ws.Range("R2").Value = 9195.845729102761
$ws.Range("S2").Value = 0.08783165599917633
$ws.Range("T2").Value = 0.08783165599917631

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.02514266666667
$ws.Range("H3").Value = 45.075428
$ws.Range("I3").Value = 0.1401726531301337
$ws.Range("J3").Value = 0.1401726531301337
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.991529999999999
$ws.Range("N3").Value = 26.97459
$ws.Range("O3").Value = 0.08284967558015671
$ws.Range("P3").Value = 0.08284967558015671
$ws.Range("Q3").Value = 135.0990210416133
$ws.Range("R3").Value = 1215.89118937452
$ws.Range("S3").Value = 0.01161325883704142
$ws.Range("T3").Value = 0.01161325883704141

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.02514266666667
$ws.Range("H4").Value = 45.075428
$ws.Range("I4").Value = 0.1401726531301337
$ws.Range("J4").Value = 0.1401726531301337
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.76843933333333
$ws.Range("N4").Value = 32.305318
$ws.Range("O4").Value = 0.09922245772090688
$ws.Range("P4").Value = 0.09922245772090688
$ws.Range("Q4").Value = 161.7973372806782
$ws.Range("R4").Value = 1456.176035526104
$ws.Range("S4").Value = 0.01390827514883204
$ws.Range("T4").Value = 0.01390827514883203

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.02514266666667
$ws.Range("H5").Value = 45.075428
$ws.Range("I5").Value = 0.1401726531301337
$ws.Range("J5").Value = 0.1401726531301337
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 20.764887
$ws.Range("N5").Value = 62.294661
$ws.Range("O5").Value = 0.1913316367079478
$ws.Range("P5").Value = 0.1913316367079478
$ws.Range("Q5").Value = 311.995389632212
$ws.Range("R5").Value = 2807.958506689908
$ws.Range("S5").Value = 0.02681946314508393
$ws.Range("T5").Value = 0.02681946314508392

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 33.26311566666666
$ws.Range("H6").Value = 99.78934699999999
$ws.Range("I6").Value = 0.3103184627135109
$ws.Range("J6").Value = 0.3103184627135109
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 68.00339
$ws.Range("N6").Value = 204.01017
$ws.Range("O6").Value = 0.6265962299909886
$ws.Range("P6").Value = 0.6265962299909885
$ws.Range("Q6").Value = 2262.004627295443
$ws.Range("R6").Value = 20358.04164565899
$ws.Range("S6").Value = 0.1944443788328851
$ws.Range("T6").Value = 0.1944443788328851

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 33.26311566666666
$ws.Range("H7").Value = 99.78934699999999
$ws.Range("I7").Value = 0.3103184627135109
$ws.Range("J7").Value = 0.3103184627135109
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.991529999999999
$ws.Range("N7").Value = 26.97459
$ws.Range("O7").Value = 0.08284967558015671
$ws.Range("P7").Value = 0.08284967558015671
$ws.Range("Q7").Value = 299.0863024103033
$ws.Range("R7").Value = 2691.77672169273
$ws.Range("S7").Value = 0.02570978396234734
$ws.Range("T7").Value = 0.02570978396234733

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 33.26311566666666
$ws.Range("H8").Value = 99.78934699999999
$ws.Range("I8").Value = 0.3103184627135109
$ws.Range("J8").Value = 0.3103184627135109
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.76843933333333
$ws.Range("N8").Value = 32.305318
$ws.Range("O8").Value = 0.09922245772090688
$ws.Range("P8").Value = 0.09922245772090688
$ws.Range("Q8").Value = 358.1918430941495
$ws.Range("R8").Value = 3223.726587847346
$ws.Range("S8").Value = 0.03079056054660816
$ws.Range("T8").Value = 0.03079056054660815

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 33.26311566666666
$ws.Range("H9").Value = 99.78934699999999
$ws.Range("I9").Value = 0.3103184627135109
$ws.Range("J9").Value = 0.3103184627135109
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 20.764887
$ws.Range("N9").Value = 62.294661
$ws.Range("O9").Value = 0.1913316367079478
$ws.Range("P9").Value = 0.1913316367079478
$ws.Range("Q9").Value = 690.7048380862628
$ws.Range("R9").Value = 6216.343542776366
$ws.Range("S9").Value = 0.05937373937167031
$ws.Range("T9").Value = 0.0593737393716703

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 50.12360066666667
$ws.Range("H10").Value = 150.370802
$ws.Range("I10").Value = 0.467613402797773
$ws.Range("J10").Value = 0.4676134027977729
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 68.00339
$ws.Range("N10").Value = 204.01017
$ws.Range("O10").Value = 0.6265962299909886
$ws.Range("P10").Value = 0.6265962299909885
$ws.Range("Q10").Value = 3408.574764339593
$ws.Range("R10").Value = 30677.17287905634
$ws.Range("S10").Value = 0.2930047952863422
$ws.Range("T10").Value = 0.293004795286342

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 50.12360066666667
$ws.Range("H11").Value = 150.370802
$ws.Range("I11").Value = 0.467613402797773
$ws.Range("J11").Value = 0.4676134027977729
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 8.991529999999999
$ws.Range("N11").Value = 26.97459
$ws.Range("O11").Value = 0.08284967558015671
$ws.Range("P11").Value = 0.08284967558015671
$ws.Range("Q11").Value = 450.6878591023533
$ws.Range("R11").Value = 4056.19073192118
$ws.Range("S11").Value = 0.03874161871872864
$ws.Range("T11").Value = 0.03874161871872862

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 50.12360066666667
$ws.Range("H12").Value = 150.370802
$ws.Range("I12").Value = 0.467613402797773
$ws.Range("J12").Value = 0.4676134027977729
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.76843933333333
$ws.Range("N12").Value = 32.305318
$ws.Range("O12").Value = 0.09922245772090688
$ws.Range("P12").Value = 0.09922245772090688
$ws.Range("Q12").Value = 539.7529529472263
$ws.Range("R12").Value = 4857.776576525036
$ws.Range("S12").Value = 0.04639775108883143
$ws.Range("T12").Value = 0.04639775108883142

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 50.12360066666667
$ws.Range("H13").Value = 150.370802
$ws.Range("I13").Value = 0.467613402797773
$ws.Range("J13").Value = 0.4676134027977729
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 20.764887
$ws.Range("N13").Value = 62.294661
$ws.Range("O13").Value = 0.1913316367079478
$ws.Range("P13").Value = 0.1913316367079478
$ws.Range("Q13").Value = 1040.810903876458
$ws.Range("R13").Value = 9367.298134888122
$ws.Range("S13").Value = 0.08946923770387076
$ws.Range("T13").Value = 0.08946923770387073

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 8.778397666666669
$ws.Range("H14").Value = 26.335193
$ws.Range("I14").Value = 0.08189548135858246
$ws.Range("J14").Value = 0.08189548135858243
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 68.00339
$ws.Range("N14").Value = 204.01017
$ws.Range("O14").Value = 0.6265962299909886
$ws.Range("P14").Value = 0.6265962299909885
$ws.Range("Q14").Value = 596.9608001014234
$ws.Range("R14").Value = 5372.64720091281
$ws.Range("S14").Value = 0.05131539987258506
$ws.Range("T14").Value = 0.05131539987258503

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 8.778397666666669
$ws.Range("H15").Value = 26.335193
$ws.Range("I15").Value = 0.08189548135858246
$ws.Range("J15").Value = 0.08189548135858243
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 8.991529999999999
$ws.Range("N15").Value = 26.97459
$ws.Range("O15").Value = 0.08284967558015671
$ws.Range("P15").Value = 0.08284967558015671
$ws.Range("Q15").Value = 78.93122597176334
$ws.Range("R15").Value = 710.3810337458701
$ws.Range("S15").Value = 0.006785014062039328
$ws.Range("T15").Value = 0.006785014062039327

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 8.778397666666669
$ws.Range("H16").Value = 26.335193
$ws.Range("I16").Value = 0.08189548135858246
$ws.Range("J16").Value = 0.08189548135858243
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 10.76843933333333
$ws.Range("N16").Value = 32.305318
$ws.Range("O16").Value = 0.09922245772090688
$ws.Range("P16").Value = 0.09922245772090688
$ws.Range("Q16").Value = 94.52964271737491
$ws.Range("R16").Value = 850.7667844563741
$ws.Range("S16").Value = 0.008125870936635266
$ws.Range("T16").Value = 0.008125870936635263

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 8.778397666666669
$ws.Range("H17").Value = 26.335193
$ws.Range("I17").Value = 0.08189548135858246
$ws.Range("J17").Value = 0.08189548135858243
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 20.764887
$ws.Range("N17").Value = 62.294661
$ws.Range("O17").Value = 0.1913316367079478
$ws.Range("P17").Value = 0.1913316367079478
$ws.Range("Q17").Value = 182.282435589397
$ws.Range("R17").Value = 1640.541920304573
$ws.Range("S17").Value = 0.01566919648732281
$ws.Range("T17").Value = 0.0156691964873228
